$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Ccl21b"
$ws.Cells.Item(2, 3).Value = "Ackr4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.134289
$ws.Cells.Item(2, 8).Value = 0.402867
$ws.Cells.Item(2, 9).Value = 0.3678949098679525
$ws.Cells.Item(2, 10).Value = 0.3678949098679525
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.2315733333333333
$ws.Cells.Item(2, 14).Value = 0.69472
$ws.Cells.Item(2, 15).Value = 0.2371078251520917
$ws.Cells.Item(2, 16).Value = 0.2371078251520917
$ws.Cells.Item(2, 17).Value = 0.03109775136
$ws.Cells.Item(2, 18).Value = 0.27987976224
$ws.Cells.Item(2, 19).Value = 0.08723076196331501
$ws.Cells.Item(2, 20).Value = 0.08723076196331501

# Row 3
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Ccl21b"
$ws.Cells.Item(3, 3).Value = "Ackr4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.134289
$ws.Cells.Item(3, 8).Value = 0.402867
$ws.Cells.Item(3, 9).Value = 0.3678949098679525
$ws.Cells.Item(3, 10).Value = 0.3678949098679525
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.7200953333333334
$ws.Cells.Item(3, 14).Value = 2.160286
$ws.Cells.Item(3, 15).Value = 0.7373052671097876
$ws.Cells.Item(3, 16).Value = 0.7373052671097876
$ws.Cells.Item(3, 17).Value = 0.09670088221800002
$ws.Cells.Item(3, 18).Value = 0.8703079399620002
$ws.Cells.Item(3, 19).Value = 0.271250854788522
$ws.Cells.Item(3, 20).Value = 0.271250854788522

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Ccl21b"
$ws.Cells.Item(4, 3).Value = "Ackr4"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.134289
$ws.Cells.Item(4, 8).Value = 0.402867
$ws.Cells.Item(4, 9).Value = 0.3678949098679525
$ws.Cells.Item(4, 10).Value = 0.3678949098679525
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.02498966666666666
$ws.Cells.Item(4, 14).Value = 0.074969
$ws.Cells.Item(4, 15).Value = 0.02558690773812063
$ws.Cells.Item(4, 16).Value = 0.02558690773812063
$ws.Cells.Item(4, 17).Value = 0.003355837347
$ws.Cells.Item(4, 18).Value = 0.030202536123
$ws.Cells.Item(4, 19).Value = 0.009413293116115504
$ws.Cells.Item(4, 20).Value = 0.009413293116115504

# Row 5
$ws.Cells.Item(5, 1).Value = "sCs"
$ws.Cells.Item(5, 2).Value = "Ccl21b"
$ws.Cells.Item(5, 3).Value = "Ackr4"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.230731
$ws.Cells.Item(5, 8).Value = 0.692193
$ws.Cells.Item(5, 9).Value = 0.6321050901320475
$ws.Cells.Item(5, 10).Value = 0.6321050901320475
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.2315733333333333
$ws.Cells.Item(5, 14).Value = 0.69472
$ws.Cells.Item(5, 15).Value = 0.2371078251520917
$ws.Cells.Item(5, 16).Value = 0.2371078251520917
$ws.Cells.Item(5, 17).Value = 0.05343114677333333
$ws.Cells.Item(5, 18).Value = 0.4808803209599999
$ws.Cells.Item(5, 19).Value = 0.1498770631887767
$ws.Cells.Item(5, 20).Value = 0.1498770631887767

# Row 6
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Ccl21b"
$ws.Cells.Item(6, 3).Value = "Ackr4"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.230731
$ws.Cells.Item(6, 8).Value = 0.692193
$ws.Cells.Item(6, 9).Value = 0.6321050901320475
$ws.Cells.Item(6, 10).Value = 0.6321050901320475
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.7200953333333334
$ws.Cells.Item(6, 14).Value = 2.160286
$ws.Cells.Item(6, 15).Value = 0.7373052671097876
$ws.Cells.Item(6, 16).Value = 0.7373052671097876
$ws.Cells.Item(6, 17).Value = 0.1661483163553333
$ws.Cells.Item(6, 18).Value = 1.495334847198
$ws.Cells.Item(6, 19).Value = 0.4660544123212657
$ws.Cells.Item(6, 20).Value = 0.4660544123212657

# Row 7
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Ccl21b"
$ws.Cells.Item(7, 3).Value = "Ackr4"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.230731
$ws.Cells.Item(7, 8).Value = 0.692193
$ws.Cells.Item(7, 9).Value = 0.6321050901320475
$ws.Cells.Item(7, 10).Value = 0.6321050901320475
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.02498966666666666
$ws.Cells.Item(7, 14).Value = 0.074969
$ws.Cells.Item(7, 15).Value = 0.02558690773812063
$ws.Cells.Item(7, 16).Value = 0.02558690773812063
$ws.Cells.Item(7, 17).Value = 0.005765890779666666
$ws.Cells.Item(7, 18).Value = 0.05189301701699999
$ws.Cells.Item(7, 19).Value = 0.01617361462200512
$ws.Cells.Item(7, 20).Value = 0.01617361462200512
